# Apply the diff: header D/E swap (time/distance), data value changes, and new rows 38-48
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1,1).Value = "alg_name"
$ws.Cells.Item(1,2).Value = "buses"
$ws.Cells.Item(1,3).Value = "passengers"
$ws.Cells.Item(1,4).Value = "time"
$ws.Cells.Item(1,5).Value = "distance"
$ws.Cells.Item(1,6).Value = "wait"

$ws.Cells.Item(2,1).Value = "new_heuristic"
$ws.Cells.Item(2,2).Value = 2
$ws.Cells.Item(2,3).Value = 2
$ws.Cells.Item(2,4).Value = 642.7679219733545
$ws.Cells.Item(2,5).Value = 20.31457473333333
$ws.Cells.Item(2,6).Value = 0

$ws.Cells.Item(3,1).Value = "greedy"
$ws.Cells.Item(3,2).Value = 2
$ws.Cells.Item(3,3).Value = 2
$ws.Cells.Item(3,4).Value = 642.7679219733545
$ws.Cells.Item(3,5).Value = 20.31457473333333
$ws.Cells.Item(3,6).Value = 0

$ws.Cells.Item(4,1).Value = "new_heuristic"
$ws.Cells.Item(4,2).Value = 2
$ws.Cells.Item(4,3).Value = 8
$ws.Cells.Item(4,4).Value = 597.75498125187
$ws.Cells.Item(4,5).Value = 30.57566683333334
$ws.Cells.Item(4,6).Value = 0

$ws.Cells.Item(5,1).Value = "greedy"
$ws.Cells.Item(5,2).Value = 2
$ws.Cells.Item(5,3).Value = 8
$ws.Cells.Item(5,4).Value = 477.0098043965991
$ws.Cells.Item(5,5).Value = 22.88711776666668
$ws.Cells.Item(5,6).Value = 0

$ws.Cells.Item(6,1).Value = "new_heuristic"
$ws.Cells.Item(6,2).Value = 2
$ws.Cells.Item(6,3).Value = 14
$ws.Cells.Item(6,4).Value = 948.1799861321342
$ws.Cells.Item(6,5).Value = 18.92675466666674
$ws.Cells.Item(6,6).Value = 0

$ws.Cells.Item(7,1).Value = "greedy"
$ws.Cells.Item(7,2).Value = 2
$ws.Cells.Item(7,3).Value = 14
$ws.Cells.Item(7,4).Value = 128.2696691792225
$ws.Cells.Item(7,5).Value = 40.58960979999995
$ws.Cells.Item(7,6).Value = 0

$ws.Cells.Item(8,1).Value = "new_heuristic"
$ws.Cells.Item(8,2).Value = 2
$ws.Cells.Item(8,3).Value = 20
$ws.Cells.Item(8,4).Value = 874.5798335453146
$ws.Cells.Item(8,5).Value = 45.60432756666671
$ws.Cells.Item(8,6).Value = 0

$ws.Cells.Item(9,1).Value = "greedy"
$ws.Cells.Item(9,2).Value = 2
$ws.Cells.Item(9,3).Value = 20
$ws.Cells.Item(9,4).Value = 898.7174345861422
$ws.Cells.Item(9,5).Value = 44.63159026666665
$ws.Cells.Item(9,6).Value = 0

$ws.Cells.Item(10,1).Value = "new_heuristic"
$ws.Cells.Item(10,2).Value = 2
$ws.Cells.Item(10,3).Value = 26
$ws.Cells.Item(10,4).Value = 428.7113177718129
$ws.Cells.Item(10,5).Value = 4.48313306666671
$ws.Cells.Item(10,6).Value = 0

$ws.Cells.Item(11,1).Value = "greedy"
$ws.Cells.Item(11,2).Value = 2
$ws.Cells.Item(11,3).Value = 26
$ws.Cells.Item(11,4).Value = 363.6739013949409
$ws.Cells.Item(11,5).Value = 28.20906950000006
$ws.Cells.Item(11,6).Value = 0

$ws.Cells.Item(12,1).Value = "new_heuristic"
$ws.Cells.Item(12,2).Value = 7
$ws.Cells.Item(12,3).Value = 7
$ws.Cells.Item(12,4).Value = 547.7573367974255
$ws.Cells.Item(12,5).Value = 17.45390368333324
$ws.Cells.Item(12,6).Value = 0

$ws.Cells.Item(13,1).Value = "greedy"
$ws.Cells.Item(13,2).Value = 7
$ws.Cells.Item(13,3).Value = 7
$ws.Cells.Item(13,4).Value = 547.7573367974255
$ws.Cells.Item(13,5).Value = 17.45390368333324
$ws.Cells.Item(13,6).Value = 0

$ws.Cells.Item(14,1).Value = "new_heuristic"
$ws.Cells.Item(14,2).Value = 7
$ws.Cells.Item(14,3).Value = 13
$ws.Cells.Item(14,4).Value = 721.0845278754132
$ws.Cells.Item(14,5).Value = 47.74245885000005
$ws.Cells.Item(14,6).Value = 0

$ws.Cells.Item(15,1).Value = "greedy"
$ws.Cells.Item(15,2).Value = 7
$ws.Cells.Item(15,3).Value = 13
$ws.Cells.Item(15,4).Value = 731.4994873754331
$ws.Cells.Item(15,5).Value = 22.02163988333336
$ws.Cells.Item(15,6).Value = 0

$ws.Cells.Item(16,1).Value = "new_heuristic"
$ws.Cells.Item(16,2).Value = 7
$ws.Cells.Item(16,3).Value = 19
$ws.Cells.Item(16,4).Value = 841.5848877596436
$ws.Cells.Item(16,5).Value = 27.62617048333323
$ws.Cells.Item(16,6).Value = 0

$ws.Cells.Item(17,1).Value = "greedy"
$ws.Cells.Item(17,2).Value = 7
$ws.Cells.Item(17,3).Value = 19
$ws.Cells.Item(17,4).Value = 293.0564817967825
$ws.Cells.Item(17,5).Value = 8.144220766666649
$ws.Cells.Item(17,6).Value = 0

$ws.Cells.Item(18,1).Value = "new_heuristic"
$ws.Cells.Item(18,2).Value = 7
$ws.Cells.Item(18,3).Value = 25
$ws.Cells.Item(18,4).Value = 848.6064657588722
$ws.Cells.Item(18,5).Value = 46.60815641666659
$ws.Cells.Item(18,6).Value = 0

$ws.Cells.Item(19,1).Value = "greedy"
$ws.Cells.Item(19,2).Value = 7
$ws.Cells.Item(19,3).Value = 25
$ws.Cells.Item(19,4).Value = 303.9292160845362
$ws.Cells.Item(19,5).Value = 53.70323998333333
$ws.Cells.Item(19,6).Value = 0

$ws.Cells.Item(20,1).Value = "new_heuristic"
$ws.Cells.Item(20,2).Value = 7
$ws.Cells.Item(20,3).Value = 31
$ws.Cells.Item(20,4).Value = 448.2121697727125
$ws.Cells.Item(20,5).Value = 11.83629648333317
$ws.Cells.Item(20,6).Value = 0

$ws.Cells.Item(21,1).Value = "greedy"
$ws.Cells.Item(21,2).Value = 7
$ws.Cells.Item(21,3).Value = 31
$ws.Cells.Item(21,4).Value = 186.7221115144202
$ws.Cells.Item(21,5).Value = 4.422495866666623
$ws.Cells.Item(21,6).Value = 0

$ws.Cells.Item(22,1).Value = "new_heuristic"
$ws.Cells.Item(22,2).Value = 12
$ws.Cells.Item(22,3).Value = 12
$ws.Cells.Item(22,4).Value = 67.15739480615593
$ws.Cells.Item(22,5).Value = 56.78359380000052
$ws.Cells.Item(22,6).Value = 0

$ws.Cells.Item(23,1).Value = "greedy"
$ws.Cells.Item(23,2).Value = 12
$ws.Cells.Item(23,3).Value = 12
$ws.Cells.Item(23,4).Value = 67.15739480615593
$ws.Cells.Item(23,5).Value = 56.78359380000052
$ws.Cells.Item(23,6).Value = 0

$ws.Cells.Item(24,1).Value = "new_heuristic"
$ws.Cells.Item(24,2).Value = 12
$ws.Cells.Item(24,3).Value = 18
$ws.Cells.Item(24,4).Value = 861.7921812601853
$ws.Cells.Item(24,5).Value = 6.86291779999965
$ws.Cells.Item(24,6).Value = 0

$ws.Cells.Item(25,1).Value = "greedy"
$ws.Cells.Item(25,2).Value = 12
$ws.Cells.Item(25,3).Value = 18
$ws.Cells.Item(25,4).Value = 57.91654975304846
$ws.Cells.Item(25,5).Value = 11.67415100000017
$ws.Cells.Item(25,6).Value = 0

$ws.Cells.Item(26,1).Value = "new_heuristic"
$ws.Cells.Item(26,2).Value = 12
$ws.Cells.Item(26,3).Value = 24
$ws.Cells.Item(26,4).Value = 868.2245132206008
$ws.Cells.Item(26,5).Value = 29.58016979999991
$ws.Cells.Item(26,6).Value = 0

$ws.Cells.Item(27,1).Value = "greedy"
$ws.Cells.Item(27,2).Value = 12
$ws.Cells.Item(27,3).Value = 24
$ws.Cells.Item(27,4).Value = 481.2080727429129
$ws.Cells.Item(27,5).Value = 12.16152280000006
$ws.Cells.Item(27,6).Value = 0

$ws.Cells.Item(28,1).Value = "greedy"
$ws.Cells.Item(28,2).Value = 12
$ws.Cells.Item(28,3).Value = 30
$ws.Cells.Item(28,4).Value = 560.4315571651096
$ws.Cells.Item(28,5).Value = 36.63915679999991
$ws.Cells.Item(28,6).Value = 0

$ws.Cells.Item(29,1).Value = "greedy"
$ws.Cells.Item(29,2).Value = 12
$ws.Cells.Item(29,3).Value = 36
$ws.Cells.Item(29,4).Value = 144.2584183252184
$ws.Cells.Item(29,5).Value = 23.99181600000065
$ws.Cells.Item(29,6).Value = 0

$ws.Cells.Item(30,1).Value = "new_heuristic"
$ws.Cells.Item(30,2).Value = 17
$ws.Cells.Item(30,3).Value = 17
$ws.Cells.Item(30,4).Value = 232.783588927472
$ws.Cells.Item(30,5).Value = 47.3355512833341
$ws.Cells.Item(30,6).Value = 0

$ws.Cells.Item(31,1).Value = "greedy"
$ws.Cells.Item(31,2).Value = 17
$ws.Cells.Item(31,3).Value = 17
$ws.Cells.Item(31,4).Value = 232.783588927472
$ws.Cells.Item(31,5).Value = 47.3355512833341
$ws.Cells.Item(31,6).Value = 0

$ws.Cells.Item(32,1).Value = "greedy"
$ws.Cells.Item(32,2).Value = 17
$ws.Cells.Item(32,3).Value = 23
$ws.Cells.Item(32,4).Value = 553.6321877402952
$ws.Cells.Item(32,5).Value = 17.04317311666637
$ws.Cells.Item(32,6).Value = 0

$ws.Cells.Item(33,1).Value = "new_heuristic"
$ws.Cells.Item(33,2).Value = 22
$ws.Cells.Item(33,3).Value = 22
$ws.Cells.Item(33,4).Value = 916.5275729757268
$ws.Cells.Item(33,5).Value = 34.67926749999924
$ws.Cells.Item(33,6).Value = 0

$ws.Cells.Item(34,1).Value = "greedy"
$ws.Cells.Item(34,2).Value = 22
$ws.Cells.Item(34,3).Value = 22
$ws.Cells.Item(34,4).Value = 916.5275729757268
$ws.Cells.Item(34,5).Value = 34.67926749999924
$ws.Cells.Item(34,6).Value = 0

$ws.Cells.Item(35,1).Value = "greedy"
$ws.Cells.Item(35,2).Value = 22
$ws.Cells.Item(35,3).Value = 28
$ws.Cells.Item(35,4).Value = 85.1124626705423
$ws.Cells.Item(35,5).Value = 48.95521896666605
$ws.Cells.Item(35,6).Value = 0

$ws.Cells.Item(36,1).Value = "greedy"
$ws.Cells.Item(36,2).Value = 22
$ws.Cells.Item(36,3).Value = 34
$ws.Cells.Item(36,4).Value = 822.6587967458181
$ws.Cells.Item(36,5).Value = 4.364584666666815
$ws.Cells.Item(36,6).Value = 0

$ws.Cells.Item(37,1).Value = "greedy"
$ws.Cells.Item(37,2).Value = 27
$ws.Cells.Item(37,3).Value = 27
$ws.Cells.Item(37,4).Value = 358.2637687600218
$ws.Cells.Item(37,5).Value = 57.25282379999999
$ws.Cells.Item(37,6).Value = 0

$ws.Cells.Item(38,1).Value = "greedy"
$ws.Cells.Item(38,2).Value = 27
$ws.Cells.Item(38,3).Value = 33
$ws.Cells.Item(38,4).Value = 876.2048105942085
$ws.Cells.Item(38,5).Value = 6.243461100000786
$ws.Cells.Item(38,6).Value = 0

$ws.Cells.Item(39,1).Value = "greedy"
$ws.Cells.Item(39,2).Value = 27
$ws.Cells.Item(39,3).Value = 45
$ws.Cells.Item(39,4).Value = 124.6155862307642
$ws.Cells.Item(39,5).Value = 7.811352450000868
$ws.Cells.Item(39,6).Value = 0

$ws.Cells.Item(40,1).Value = "greedy"
$ws.Cells.Item(40,2).Value = 32
$ws.Cells.Item(40,3).Value = 38
$ws.Cells.Item(40,4).Value = 150.6909244614653
$ws.Cells.Item(40,5).Value = 7.829381866668882
$ws.Cells.Item(40,6).Value = 0

$ws.Cells.Item(41,1).Value = "greedy"
$ws.Cells.Item(41,2).Value = 32
$ws.Cells.Item(41,3).Value = 44
$ws.Cells.Item(41,4).Value = 520.1484631183557
$ws.Cells.Item(41,5).Value = 54.7035493333351
$ws.Cells.Item(41,6).Value = 0

$ws.Cells.Item(42,1).Value = "new_heuristic"
$ws.Cells.Item(42,2).Value = 37
$ws.Cells.Item(42,3).Value = 37
$ws.Cells.Item(42,4).Value = 642.058284324361
$ws.Cells.Item(42,5).Value = 40.4405168000003
$ws.Cells.Item(42,6).Value = 0

$ws.Cells.Item(43,1).Value = "greedy"
$ws.Cells.Item(43,2).Value = 37
$ws.Cells.Item(43,3).Value = 37
$ws.Cells.Item(43,4).Value = 642.058284324361
$ws.Cells.Item(43,5).Value = 40.4405168000003
$ws.Cells.Item(43,6).Value = 0

$ws.Cells.Item(44,1).Value = "greedy"
$ws.Cells.Item(44,2).Value = 37
$ws.Cells.Item(44,3).Value = 43
$ws.Cells.Item(44,4).Value = 836.0277959888335
$ws.Cells.Item(44,5).Value = 40.34617185000207
$ws.Cells.Item(44,6).Value = 0

$ws.Cells.Item(45,1).Value = "greedy"
$ws.Cells.Item(45,2).Value = 37
$ws.Cells.Item(45,3).Value = 49
$ws.Cells.Item(45,4).Value = 524.5183579868171
$ws.Cells.Item(45,5).Value = 37.54547608333087
$ws.Cells.Item(45,6).Value = 0

$ws.Cells.Item(46,1).Value = "greedy"
$ws.Cells.Item(46,2).Value = 37
$ws.Cells.Item(46,3).Value = 55
$ws.Cells.Item(46,4).Value = 534.6393749127164
$ws.Cells.Item(46,5).Value = 43.43480576667298
$ws.Cells.Item(46,6).Value = 0

$ws.Cells.Item(47,1).Value = "new_heuristic"
$ws.Cells.Item(47,2).Value = 42
$ws.Cells.Item(47,3).Value = 42
$ws.Cells.Item(47,4).Value = 687.5053024478257
$ws.Cells.Item(47,5).Value = 40.39948190000041
$ws.Cells.Item(47,6).Value = 0

$ws.Cells.Item(48,1).Value = "greedy"
$ws.Cells.Item(48,2).Value = 42
$ws.Cells.Item(48,3).Value = 42
$ws.Cells.Item(48,4).Value = 687.5053024478257
$ws.Cells.Item(48,5).Value = 40.39948190000041
$ws.Cells.Item(48,6).Value = 0

[void]$ws.Range("I14").Select()
"edit complete"